$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# content edits (order matters for shared-string table ordering)
$ws.Range("E11").Value = "YES"
$ws.Range("C4").Value = "Grey"
$ws.Range("C7").Value = "Grey"

# new rows
$ws.Range("A13").Value = 2569
$ws.Range("B13").Value = 4251564
$ws.Range("C13").Value = "Grey"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "NO"

$ws.Range("A14").Value = 6111
$ws.Range("B14").Value = 611101
$ws.Range("C14").Value = "White"
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = "NO"

$ws.Range("A15").Value = 330701
$ws.Range("B15").Value = 3307
$ws.Range("C15").Value = "White"
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = "NO"

# alignment: header row vertical center (already horizontal=center)
$ws.Range("A1:E1").VerticalAlignment = -4108

# alignment: columns A:B rows 2-12 stay right aligned, add vertical center
$ws.Range("A2:B12").VerticalAlignment = -4108

# alignment: columns A:B rows 13-15 -> vertical center only (general horizontal)
$ws.Range("A13:B15").VerticalAlignment = -4108

# alignment: columns C:D:E rows 2-15 -> center horizontal + vertical center
$ws.Range("C2:E15").VerticalAlignment = -4108
$ws.Range("C2:E15").HorizontalAlignment = -4108

# update selection
[void]$ws.Range("D20").Select()
